$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# rech row (row 2): faclbnd (I2) and facubnd (J2)
$ws.Range("I2").Value = 0.01
$ws.Range("J2").Value = 100

# hk row (row 3): faclbnd (I3) and facubnd (J3)
$ws.Range("I3").Value = 0.001
$ws.Range("J3").Value = 10000

# cdrn row (row 24): parubnd (H24)
$ws.Range("H24").Value = 10

# cghb row (row 27): parubnd (H27)
$ws.Range("H27").Value = 10

# update the active selection to match the new state
$null = $ws.Range("G19").Select()
